$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the 7 pen-group header rows with descriptive text instead of
# bare numbers (9-15 -> "all_in_all_out 9".."all_in_all_out 15").
$ws.Range("A5").Value = "all_in_all_out 9"
$ws.Range("A8").Value = "all_in_all_out 10"
$ws.Range("A11").Value = "all_in_all_out 11"
$ws.Range("A14").Value = "all_in_all_out 12"
$ws.Range("A17").Value = "all_in_all_out 13"
$ws.Range("A20").Value = "all_in_all_out 14"
$ws.Range("A23").Value = "all_in_all_out 15"

# Widen column A so the longer labels fit.
$ws.Columns("A").ColumnWidth = 18.33

# Leave the selection on A25, matching where editing ended up.
$ws.Range("A25").Select()
